$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "last refreshed" timestamp banner (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 29 de Mayo de 2020 a las 21:40"

# --- Estados Unidos (row 4): refreshed COVID totals ---
$ws.Range("B4").Value = 1784223
$ws.Range("C4").Value = 15762
$ws.Range("E4").Value = 1177909
$ws.Range("G4").Value = 721
$ws.Range("H4").Value = 104051

# --- Costa Rica (row 116): refreshed COVID totals ---
$ws.Range("B116").Value = 1022
$ws.Range("C116").Value = 22
$ws.Range("D116").Value = 653
$ws.Range("E116").Value = 359

# --- Mauritania's case count overtakes Togo/Cabo Verde/Ruanda, so the
#     country list (sorted descending by "Casos totales") is re-sorted:
#     rows 143-146 shift down to make room and Mauritania's updated
#     figures land in row 143, right after Guayana Francesa. ---
$ws.Range("A143").Value = "Mauritania"
$ws.Range("B143").Value = 423
$ws.Range("C143").Value = 77
$ws.Range("D143").Value = 21
$ws.Range("E143").Value = 382
$ws.Range("F143").Value = 0
$ws.Range("G143").Value = 1
$ws.Range("H143").Value = 20

$ws.Range("A144").Value = "Togo"
$ws.Range("B144").Value = 422
$ws.Range("C144").Value = 0
$ws.Range("D144").Value = 197
$ws.Range("E144").Value = 212
$ws.Range("F144").Value = 0
$ws.Range("G144").Value = 0
$ws.Range("H144").Value = 13

$ws.Range("A145").Value = "Cabo Verde"
$ws.Range("B145").Value = 390
$ws.Range("C145").Value = 0
$ws.Range("D145").Value = 155
$ws.Range("E145").Value = 231
$ws.Range("F145").Value = 0
$ws.Range("G145").Value = 0
$ws.Range("H145").Value = 4

$ws.Range("A146").Value = "Ruanda"
$ws.Range("B146").Value = 349
$ws.Range("C146").Value = 0
$ws.Range("D146").Value = 245
$ws.Range("E146").Value = 104
$ws.Range("F146").Value = 0
$ws.Range("G146").Value = 0
$ws.Range("H146").Value = 0
